# Generate Report for Handback
# Fills in the handback-result columns (Latest Target File / Latest Handback
# File / Latest Handback DateTime / Error Detail) for the 743f741e... row on
# both the zh-cn and de-de sheets, now that a (stale) handback was processed
# for that item, and widens the "Error Detail" column so the message is
# readable.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a1456284c4548b7c081ba03486618032fcd76bb5/e2e/743f741e-0d05-4bbb-9753-8b3eec852616.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f16261dd093b0e69b7f8c007adec7e27391fc9a7/e2e/743f741e-0d05-4bbb-9753-8b3eec852616.md."

$locales = @(
    @{
        sheet = "zh-cn"
        targetFileName = "743f741e-0d05-4bbb-9753-8b3eec852616.5bdb2170776f8fac65ec5384e3eb612bc42093c6.zh-cn.xlf"
        handbackDateTime = "2016-08-30 09:18:34"
        hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ba6c31aa634e8178f7b91148ba46e8981cf379a8/e2e/743f741e-0d05-4bbb-9753-8b3eec852616.md"
    },
    @{
        sheet = "de-de"
        targetFileName = "743f741e-0d05-4bbb-9753-8b3eec852616.5bdb2170776f8fac65ec5384e3eb612bc42093c6.de-de.xlf"
        handbackDateTime = "2016-08-30 09:18:41"
        hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b490ab855f9c8d94b6c4453b1987f087a226bc92/e2e/743f741e-0d05-4bbb-9753-8b3eec852616.md"
    }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.sheet)

    # I6: Latest Target File - gains the handed-back file's display name and
    # a link into the locale-specific repo (mirrors rows 2-5's I-column links).
    $ws.Hyperlinks.Add($ws.Range("I6"), $loc.hyperlinkUrl, "", "", "743f741e-0d05-4bbb-9753-8b3eec852616.md") | Out-Null

    # J6: Latest Handback File - the handed-back xlf for this locale (same
    # file name already shown in G6, "Latest Handoff File").
    $ws.Range("J6").Value = $loc.targetFileName

    # K6: Latest Handback DateTime - timestamp of this handback attempt.
    $ws.Range("K6").Value = $loc.handbackDateTime

    # P6: Error Detail - the handback was against a stale handoff version.
    $ws.Range("P6").Value = $errorDetail

    # Error Detail column needed to grow to fit the long message.
    $ws.Columns.Item(16).ColumnWidth = 39.17
}

Write-Output "Handback report generated."
